# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.215.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.582.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.39%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("E8").Value = "  -1.49%  "

# Row 9
$ws.Range("E9").Value = "  -0.51%  "

# Row 10
$ws.Range("E10").Value = "  -0.98%  "

# Row 11
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.805.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.586.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("E14").Value = "  -0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.218.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.78%  "

# Row 18
$ws.Range("E18").Value = "  -0.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "

# Row 22
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("E23").Value = "  -2.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.95%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26
$ws.Range("E26").Value = "  -0.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "

# Row 28
$ws.Range("E28").Value = "  -1.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "

# Row 30
$ws.Range("E30").Value = "  -1.41%  "

# Row 31
$ws.Range("E31").Value = "  -0.77%  "

# Row 32
$ws.Range("E32").Value = "  -1.33%  "

# Row 33
$ws.Range("E33").Value = "  -1.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.282.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.64%  "

# Row 36
$ws.Range("E36").Value = "  -0.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.604"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "

# Row 38
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("E39").Value = "  -1.74%  "

# Row 40
$ws.Range("E40").Value = "  -1.93%  "

# Row 41
$ws.Range("E41").Value = "  +3.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.770"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("E43").Value = "  -2.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.718.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48
$ws.Range("E48").Value = "  -0.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.37%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "

